$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Seed rows 86 and 88 with the "even band" formatting (border,
#    no fill, center/left alignment) by copying row 84's formats,
#    and row 87 with the "odd band" formatting (border, grey fill)
#    by copying row 3's formats.
# ---------------------------------------------------------------
$ws.Range("A84:AK84").Copy() | Out-Null
$ws.Range("A86:AK86").PasteSpecial(-4122) | Out-Null

$ws.Range("A3:AK3").Copy() | Out-Null
$ws.Range("A87:AK87").PasteSpecial(-4122) | Out-Null

$ws.Range("A84:AK84").Copy() | Out-Null
$ws.Range("A88:AK88").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------
# 2. The "work description" (P) and "work content" (AC) columns use
#    a wrap-text variant of the banded style. Row 85's P/AC cells
#    also switch to the wrap-text variant in this edit.
# ---------------------------------------------------------------
$ws.Range("P85").WrapText = $true
$ws.Range("AC85").WrapText = $true
$ws.Range("P86").WrapText = $true
$ws.Range("AC86").WrapText = $true
$ws.Range("P87").WrapText = $true
$ws.Range("AC87").WrapText = $true

# Row 86 values
$ws.Range("A86").Value = 84
$ws.Range("B86").Value = '維修'
$ws.Range("C86").Value = 2025062222
$ws.Range("D86").Value = 'E2759114061701'
$ws.Range("E86").Value = '一般件'
$ws.Range("F86").Value = 2759
$ws.Range("G86").Value = '三重仁化店'
$ws.Range("H86").Value = '新北市三重區'
$ws.Range("I86").Value = '2025-06-17 09:39:38'
$ws.Range("J86").Value = '星期二'
$ws.Range("K86").Value = '上午'
$ws.Range("L86").Value = 'HL58'
$ws.Range("M86").Value = 'HL-LIFE-ET主機'
$ws.Range("N86").Value = 5804
$ws.Range("O86").Value = '無法連線'
$ws.Range("P86").Value = '門市反應mmk4代機無法連線畫面顯示尚未連接到網路，門市重開機多次仍未跳出網路設定介面，ping60....須請台芝到店協助(LifeET無法連線網路，已嘗試重新開機多次)'
$ws.Range("Q86").Value = 'THILF02759'
$ws.Range("R86").Value = '新北一'
$ws.Range("S86").Value = '吳宗鴻'
$ws.Range("T86").Value = 1
$ws.Range("U86").Value = '已完工'
$ws.Range("V86").Value = '2025-06-17 09:57:22'
$ws.Range("W86").Value = '2025-06-17 14:02:00'
$ws.Range("X86").Value = '2025-06-17 14:32:00'
$ws.Range("Y86").Value = '2025-06-18 13:57:00'
$ws.Range("Z86").Value = 0.5
$ws.Range("AB86").Value = '到場處理'
$ws.Range("AC86").Value = '重新啟動網路連線及設定皆為正常'
$ws.Range("AK86").Value = 'O'
# Row 87 values
$ws.Range("A87").Value = 85
$ws.Range("B87").Value = '維修'
$ws.Range("C87").Value = 2025062232
$ws.Range("D87").Value = '12399114061701'
$ws.Range("E87").Value = '一般件'
$ws.Range("F87").Value = 2399
$ws.Range("G87").Value = '三重三陽店'
$ws.Range("H87").Value = '新北市三重區'
$ws.Range("I87").Value = '2025-06-17 10:21:14'
$ws.Range("J87").Value = '星期二'
$ws.Range("K87").Value = '上午'
$ws.Range("L87").Value = 'HL58'
$ws.Range("M87").Value = 'HL-LIFE-ET主機'
$ws.Range("N87").Value = 5804
$ws.Range("O87").Value = '無法連線'
$ws.Range("P87").Value = '門市反應mmk4代機無法連線畫面顯示修正連線問題，門市重開機多次仍未跳出網路設定介面，ping60不通....須請台芝到店協助'
$ws.Range("Q87").Value = 'THILF02399'
$ws.Range("R87").Value = '新北一'
$ws.Range("S87").Value = '吳宗鴻'
$ws.Range("T87").Value = 1
$ws.Range("U87").Value = '已完工'
$ws.Range("V87").Value = '2025-06-17 10:23:16'
$ws.Range("W87").Value = '2025-06-17 14:36:00'
$ws.Range("X87").Value = '2025-06-17 15:06:00'
$ws.Range("Y87").Value = '2025-06-18 14:23:00'
$ws.Range("Z87").Value = 0.5
$ws.Range("AB87").Value = '到場處理'
$ws.Range("AC87").Value = '重新啟動網路連線及設定皆為正常'
$ws.Range("AK87").Value = 'O'
# Row 88 values
$ws.Range("A88").Value = 86
$ws.Range("B88").Value = '維修'
$ws.Range("C88").Value = 2025062243
$ws.Range("D88").Value = '13601114061701'
$ws.Range("E88").Value = '一般件'
$ws.Range("F88").Value = 3601
$ws.Range("G88").Value = '北縣重富店'
$ws.Range("H88").Value = '新北市三重區'
$ws.Range("I88").Value = '2025-06-17 11:25:37'
$ws.Range("J88").Value = '星期二'
$ws.Range("K88").Value = '上午'
$ws.Range("L88").Value = 'HL58'
$ws.Range("M88").Value = 'HL-LIFE-ET主機'
$ws.Range("N88").Value = 5804
$ws.Range("O88").Value = '無法連線'
$ws.Range("P88").Value = '門市反應MMK四代機無法連線，門市重開機多次仍未跳出網路設定介面，ping60不通....須請台芝到店協助'
$ws.Range("Q88").Value = 'THILF03601'
$ws.Range("R88").Value = '新北一'
$ws.Range("S88").Value = '吳宗鴻'
$ws.Range("T88").Value = 1
$ws.Range("U88").Value = '已完工'
$ws.Range("V88").Value = '2025-06-17 11:27:01'
$ws.Range("W88").Value = '2025-06-17 13:49:00'
$ws.Range("X88").Value = '2025-06-17 14:19:00'
$ws.Range("Y88").Value = '2025-06-18 15:27:00'
$ws.Range("Z88").Value = 0.5
$ws.Range("AB88").Value = '到場處理'
$ws.Range("AC88").Value = 'Switch7號孔8號孔接反'
$ws.Range("AK88").Value = 'O'


# ---------------------------------------------------------------
# 3. Update the printable area to include the three new rows and
#    move the active-cell selection the way the source file has it.
# ---------------------------------------------------------------
$ws.PageSetup.PrintArea = "'Report'!`$A`$1:`$AK`$88"
$ws.Range("AC85").Select() | Out-Null
